# Handback status report refresh.
#
# The localization pipeline re-ran the handoff/handback cycle for the
# "c021d156-...md" source file, producing a new xliff round-trip. That
# bumps the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for that file's row on both the "zh-cn" and "de-de" sheets,
# and the Overview sheet's "Latest HO Xliff Generate Date" (which mirrors
# the de-de handoff time) for that same file follows suit.
#
# All other cells (file names, paths, statuses, hyperlinks, etc.) are
# untouched.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

# zh-cn sheet, row 2 (c021d156-...md): Correspond Handoff / Handback Datetime
$zhcn.Range("H2").Value = "2016-08-16 06:42:25"
$zhcn.Range("K2").Value = "2016-08-16 06:42:41"

# de-de sheet, row 2 (c021d156-...md): Correspond Handoff / Handback Datetime
$dede.Range("H2").Value = "2016-08-16 06:42:30"
$dede.Range("K2").Value = "2016-08-16 06:42:47"

# Overview sheet, row 2 (c021d156-...md): Latest HO Xliff Generate Date
$overview.Range("G2").Value = "2016-08-16 06:42:30"
